# Fruta / hortaliza, semanal
# Insert a new daily price record at row 4 (pushing existing rows down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44956
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100108
$ws.Range("H4").Value = "Tropicales y subtropicales"
$ws.Range("I4").Value = 100108002
$ws.Range("J4").Value = "Mango"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6750
$ws.Range("Q4").Value = "$/bandeja 4 kilos"
$ws.Range("R4").Value = "Perú"
$ws.Range("S4").Value = 1688
$ws.Range("T4").Value = 4
